$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1945.8334
$ws.Range("J17").Value = 1986.3636
$ws.Range("L17").Value = 5959.0908
$ws.Range("N17").Value = -6295.0908

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 25888
$ws.Range("J44").Value = 25888
$ws.Range("L44").Value = 25888
$ws.Range("N44").Value = -26812

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 177.57143
$ws.Range("I107").Value = 177.57143
$ws.Range("K107").Value = 177.57143
$ws.Range("M107").Value = 1742.42857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1296.742
$ws.Range("I132").Value = 1010.5
$ws.Range("K132").Value = 3031.5
$ws.Range("M132").Value = -501.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2544.2266
$ws.Range("I137").Value = 1541.2069
$ws.Range("K137").Value = 4623.620699999999
$ws.Range("M137").Value = -2073.620699999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2464.07
$ws.Range("I138").Value = 823.4286
$ws.Range("K138").Value = 2470.2858
$ws.Range("M138").Value = 2669.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7251.4
$ws.Range("I2").Value = 438.27274
$ws.Range("K2").Value = 438.27274
$ws.Range("M2").Value = -325.27274

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4278.091
$ws.Range("I32").Value = 2826.12
$ws.Range("K32").Value = 2826.12
$ws.Range("M32").Value = -2539.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4119.3096
$ws.Range("I61").Value = 3179.4358
$ws.Range("J61").Value = 16337.667
$ws.Range("K61").Value = 3179.4358
$ws.Range("L61").Value = 16337.667
$ws.Range("M61").Value = -2967.4358
$ws.Range("N61").Value = -16761.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1863.8
$ws.Range("I97").Value = 1169.6666
$ws.Range("J97").Value = 2905
$ws.Range("K97").Value = 1169.6666
$ws.Range("L97").Value = 2905
$ws.Range("M97").Value = -673.6666
$ws.Range("N97").Value = -3897

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3043.5
$ws.Range("I110").Value = 2409.6206
$ws.Range("J110").Value = 9171
$ws.Range("K110").Value = 2409.6206
$ws.Range("L110").Value = 9171
$ws.Range("M110").Value = -364.6206000000002
$ws.Range("N110").Value = -13261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 7251.4
$ws.Range("I116").Value = 438.27274
$ws.Range("K116").Value = 438.27274
$ws.Range("M116").Value = 1855.72726

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4119.3096
$ws.Range("I136").Value = 3179.4358
$ws.Range("J136").Value = 16337.667
$ws.Range("K136").Value = 9538.307400000002
$ws.Range("L136").Value = 49013.001
$ws.Range("M136").Value = -6988.307400000002
$ws.Range("N136").Value = -54113.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7251.4
$ws.Range("I3").Value = 438.27274
$ws.Range("K3").Value = 438.27274
$ws.Range("M3").Value = -324.27274

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 851.9091
$ws.Range("J80").Value = 833.375
$ws.Range("L80").Value = 833.375
$ws.Range("N80").Value = -2829.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 851.9091
$ws.Range("J83").Value = 833.375
$ws.Range("L83").Value = 4166.875
$ws.Range("N83").Value = -14150.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2929.611
$ws.Range("J94").Value = 3941.8
$ws.Range("L94").Value = 3941.8
$ws.Range("N94").Value = -4843.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6125.55
$ws.Range("I105").Value = 1688.75
$ws.Range("J105").Value = 12780.75
$ws.Range("K105").Value = 1688.75
$ws.Range("L105").Value = 12780.75
$ws.Range("M105").Value = 58.25
$ws.Range("N105").Value = -16274.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2786.8333
$ws.Range("J134").Value = 849.5
$ws.Range("L134").Value = 2548.5
$ws.Range("N134").Value = -7618.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6145.2856
$ws.Range("I62").Value = 3459.9285
$ws.Range("J62").Value = 11516
$ws.Range("K62").Value = 3459.9285
$ws.Range("L62").Value = 11516
$ws.Range("M62").Value = -2835.9285
$ws.Range("N62").Value = -12764

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 6145.2856
$ws.Range("I65").Value = 3459.9285
$ws.Range("J65").Value = 11516
$ws.Range("K65").Value = 17299.6425
$ws.Range("L65").Value = 57580
$ws.Range("M65").Value = -14179.6425
$ws.Range("N65").Value = -63820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3886.9565
$ws.Range("I132").Value = 3107.7273
$ws.Range("K132").Value = 9323.1819
$ws.Range("M132").Value = -6793.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 107.333336
$ws.Range("I44").Value = 52.2
$ws.Range("J44").Value = 176.25
$ws.Range("K44").Value = 156.6
$ws.Range("L44").Value = 528.75
$ws.Range("M44").Value = 241.4
$ws.Range("N44").Value = -1324.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4170.4
$ws.Range("I80").Value = 1900
$ws.Range("J80").Value = 5143.4287
$ws.Range("K80").Value = 5700
$ws.Range("L80").Value = 15430.2861
$ws.Range("M80").Value = -4764
$ws.Range("N80").Value = -17302.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 4170.4
$ws.Range("I83").Value = 1900
$ws.Range("J83").Value = 5143.4287
$ws.Range("K83").Value = 17100
$ws.Range("L83").Value = 46290.85830000001
$ws.Range("M83").Value = -12420
$ws.Range("N83").Value = -55650.85830000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15213423
$ws.Range("J131").Value = 9724094
$ws.Range("L131").Value = 29172282
$ws.Range("N131").Value = -29182362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 83972.75
$ws.Range("I137").Value = 697.4545000000001
$ws.Range("J137").Value = 1000001
$ws.Range("K137").Value = 2092.3635
$ws.Range("L137").Value = 3000003
$ws.Range("M137").Value = 3007.6365
$ws.Range("N137").Value = -3010203

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2778.6365
$ws.Range("I140").Value = 2397.1
$ws.Range("K140").Value = 7191.299999999999
$ws.Range("M140").Value = -2011.299999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 1500
$ws.Range("J3").Value = 1500
$ws.Range("L3").Value = 1500
$ws.Range("N3").Value = -1724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H15").Value = 1500
$ws.Range("J15").Value = 1500
$ws.Range("L15").Value = 1500
$ws.Range("N15").Value = -1840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 1500
$ws.Range("J25").Value = 1500
$ws.Range("L25").Value = 1500
$ws.Range("N25").Value = -1960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3711.9211
$ws.Range("I46").Value = 2035.8572
$ws.Range("K46").Value = 2035.8572
$ws.Range("M46").Value = -1847.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1784.9259
$ws.Range("I55").Value = 520.9286
$ws.Range("J55").Value = 3146.1538
$ws.Range("K55").Value = 520.9286
$ws.Range("L55").Value = 3146.1538
$ws.Range("M55").Value = -347.9286
$ws.Range("N55").Value = -3492.1538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3884.0571
$ws.Range("I61").Value = 3194.111
$ws.Range("K61").Value = 3194.111
$ws.Range("M61").Value = -2992.111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3884.0571
$ws.Range("I113").Value = 3194.111
$ws.Range("K113").Value = 3194.111
$ws.Range("M113").Value = -1024.111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3860.561
$ws.Range("I132").Value = 2978.4546
$ws.Range("J132").Value = 4881.9473
$ws.Range("K132").Value = 8935.363799999999
$ws.Range("L132").Value = 14645.8419
$ws.Range("M132").Value = -6405.363799999999
$ws.Range("N132").Value = -19705.8419

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1129.4783
$ws.Range("J113").Value = 1049.75
$ws.Range("L113").Value = 3149.25
$ws.Range("N113").Value = -7489.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1592.8334
$ws.Range("I132").Value = 801.3333
$ws.Range("K132").Value = 2403.9999
$ws.Range("M132").Value = 126.0001000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3723.0527
$ws.Range("I136").Value = 1983.375
$ws.Range("J136").Value = 13001.333
$ws.Range("K136").Value = 5950.125
$ws.Range("L136").Value = 39003.999
$ws.Range("M136").Value = -3400.125
$ws.Range("N136").Value = -44103.999
